# Updated symbol list on Tue Dec 13 04:42:55 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Price (column D) updates
Set-TextValue "D2" "268.56"
Set-TextValue "D4" "6.247"
Set-TextValue "D5" "0.06207"
Set-TextValue "D6" "3.568"
Set-TextValue "D7" "6.539"
Set-TextValue "D8" "1.389"
Set-TextValue "D9" "0.8259"
Set-TextValue "D10" "0.1639"
Set-TextValue "D11" "0.08273"
Set-TextValue "D12" "0.03563"
Set-TextValue "D13" "0.03186"
Set-TextValue "D14" "0.09197"
Set-TextValue "D15" "3.762"
Set-TextValue "D16" "0.001628"
Set-TextValue "D17" "0.04669"
Set-TextValue "D18" "0.006443"
Set-TextValue "D19" "0.006208"
Set-TextValue "D20" "0.001068"
Set-TextValue "D23" "2.295"
Set-TextValue "D24" "0.01368"
Set-TextValue "D25" "0.3289"
Set-TextValue "D28" "0.0002712"
Set-TextValue "D40" "0.04719"
Set-TextValue "D41" "0.006965"

# Row 42/43: coin symbol ordering swapped (BKEXToken now ranked above CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1122"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003461"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.01158"
Set-TextValue "D45" "0.00006356"
Set-TextValue "D46" "0.0009898"

Set-TextValue "D48" "0.8023"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

Set-TextValue "D49" "0.002019"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
